# CleartripTestData.xlsx - "test data changes"
#
# 1. Fix the typo'd "From date" in HotelSearch!B2 ("15/0/2020" -> "15/09/2020").
#    ("To date" in C2, "02/10/2020", is left as-is.)
# 2. Move the sheet selection on HotelSearch from B7 to C7.
# 3. Grow the saved workbook window height (cosmetic view metadata).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HotelSearch")

# --- 1. Correct the malformed date in the "From date" column ---
$ws.Range("B2").Value = "15/09/2020"

# --- 2. Update the active selection/cell ---
$ws.Activate() | Out-Null
$ws.Range("C7").Select() | Out-Null

# --- 3. Resize the workbook window (windowHeight 3885 -> 4785) ---
$excel.ActiveWindow.Height = 4785
